# Applies the cryptos-list price/volume refresh (plus a few coin
# re-rankings / swaps) described by the commit:
#   "Updated cryptos list on Sat Jan  6 08:46:28 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every B/C/D/E cell in this sheet holds text (percentages such as
# "  -0.07%  " and prices kept as literal strings like "43.996.93" or
# "94.60"), so each write below is forced into text mode with a leading
# apostrophe - exactly like typing '43.96 into a cell - so Excel does not
# silently coerce it into a Double and lose formatting (trailing zeros,
# thousand-dot groupings, the padding spaces around the percentages, etc).
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
}

Set-TextValue $ws.Range("D2") "43.954.46"
Set-TextValue $ws.Range("E2") "  -0.07%  "

Set-TextValue $ws.Range("D3") "2.232.01"
Set-TextValue $ws.Range("E3") "  -0.98%  "

Set-TextValue $ws.Range("D4") "1.01"
Set-TextValue $ws.Range("E4") "  +0.33%  "

Set-TextValue $ws.Range("D5") "304.42"
Set-TextValue $ws.Range("E5") "  -4.70%  "

Set-TextValue $ws.Range("D6") "94.60"
Set-TextValue $ws.Range("E6") "  -6.95%  "

Set-TextValue $ws.Range("D7") "0.567"
Set-TextValue $ws.Range("E7") "  -1.64%  "

Set-TextValue $ws.Range("E8") "  +0.21%  "

Set-TextValue $ws.Range("D9") "0.521"
Set-TextValue $ws.Range("E9") "  -5.71%  "

Set-TextValue $ws.Range("D10") "34.36"
Set-TextValue $ws.Range("E10") "  -7.94%  "

Set-TextValue $ws.Range("D11") "0.0803"
Set-TextValue $ws.Range("E11") "  -3.31%  "

Set-TextValue $ws.Range("D12") "7.15"
Set-TextValue $ws.Range("E12") "  -6.00%  "

Set-TextValue $ws.Range("E13") "  -2.78%  "

Set-TextValue $ws.Range("D14") "2.574.25"
Set-TextValue $ws.Range("E14") "  -0.96%  "

Set-TextValue $ws.Range("D15") "2.235.15"
Set-TextValue $ws.Range("E15") "  -1.10%  "

Set-TextValue $ws.Range("D16") "0.813"
Set-TextValue $ws.Range("E16") "  -5.06%  "

Set-TextValue $ws.Range("D17") "13.43"
Set-TextValue $ws.Range("E17") "  -6.95%  "

Set-TextValue $ws.Range("D18") "43.840.56"
Set-TextValue $ws.Range("E18") "  -0.12%  "

Set-TextValue $ws.Range("D19") "0.0₃0951"
Set-TextValue $ws.Range("E19") "  -3.37%  "

Set-TextValue $ws.Range("D20") "12.12"
Set-TextValue $ws.Range("E20") "  -8.97%  "

Set-TextValue $ws.Range("D21") "6.14"

Set-TextValue $ws.Range("D22") "64.55"

Set-TextValue $ws.Range("D23") "236.56"
Set-TextValue $ws.Range("E23") "  +0.51%  "

Set-TextValue $ws.Range("D24") "2.89"
Set-TextValue $ws.Range("E24") "  -8.43%  "

Set-TextValue $ws.Range("E25") "  +0.51%  "

Set-TextValue $ws.Range("D26") "1.93"
Set-TextValue $ws.Range("E26") "  -8.14%  "

Set-TextValue $ws.Range("D27") "9.80"
Set-TextValue $ws.Range("E27") "  -4.06%  "

Set-TextValue $ws.Range("D28") "2.12"
Set-TextValue $ws.Range("E28") "  -3.17%  "

Set-TextValue $ws.Range("D29") "36.41"
Set-TextValue $ws.Range("E29") "  -3.50%  "

Set-TextValue $ws.Range("D30") "20.03"
Set-TextValue $ws.Range("E30") "  -0.86%  "

Set-TextValue $ws.Range("D31") "5.83"
Set-TextValue $ws.Range("E31") "  -5.66%  "

Set-TextValue $ws.Range("D32") "153.48"

Set-TextValue $ws.Range("D33") "0.0801"
Set-TextValue $ws.Range("E33") "  -5.78%  "

Set-TextValue $ws.Range("B34") "LidoDAOToken"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D34") "3.25"
Set-TextValue $ws.Range("E34") "  +7.61%  "

Set-TextValue $ws.Range("B35") "WEMIXToken"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D35") "2.61"
Set-TextValue $ws.Range("E35") "  -2.35%  "

Set-TextValue $ws.Range("D36") "0.108"
Set-TextValue $ws.Range("E36") "  -6.23%  "

Set-TextValue $ws.Range("E37") "  -0.87%  "

Set-TextValue $ws.Range("E38") "  -8.96%  "

Set-TextValue $ws.Range("D39") "14.77"
Set-TextValue $ws.Range("E39") "  -11.44%  "

Set-TextValue $ws.Range("D40") "3.34"
Set-TextValue $ws.Range("E40") "  -10.19%  "

Set-TextValue $ws.Range("D41") "3.75"
Set-TextValue $ws.Range("E41") "  -10.82%  "

Set-TextValue $ws.Range("D42") "0.0297"
Set-TextValue $ws.Range("E42") "  -5.53%  "

Set-TextValue $ws.Range("E43") "  +0.26%  "

Set-TextValue $ws.Range("D44") "1.733.11"
Set-TextValue $ws.Range("E44") "  -3.44%  "

Set-TextValue $ws.Range("D45") "83.72"
Set-TextValue $ws.Range("E45") "  +1.36%  "

Set-TextValue $ws.Range("E46") "  -7.03%  "

Set-TextValue $ws.Range("D47") "99.22"
Set-TextValue $ws.Range("E47") "  -5.20%  "

Set-TextValue $ws.Range("E48") "  -6.38%  "

Set-TextValue $ws.Range("B49") "FraxShare"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D49") "8.02"
Set-TextValue $ws.Range("E49") "  -4.06%  "

Set-TextValue $ws.Range("B50") "ordi"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue $ws.Range("D50") "68.11"
Set-TextValue $ws.Range("E50") "  -9.73%  "

Set-TextValue $ws.Range("B51") "MultiversX"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D51") "53.72"
Set-TextValue $ws.Range("E51") "  -7.48%  "
